$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '26.740.92'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'" + '  -1.30%  '
$ws.Range('E2').Style = "Normal"

# Row 3
$ws.Range('D3').Value = "'" + '1.798.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'" + '  -1.19%  '
$ws.Range('E3').Style = "Normal"

# Row 4
$ws.Range('D4').Value = "'" + '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = "'" + '  -0.11%  '
$ws.Range('E4').Style = "Normal"

# Row 5
$ws.Range('D5').Value = "'" + '309.76'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'" + '  -0.28%  '
$ws.Range('E5').Style = "Normal"

# Row 6
$ws.Range('D6').Value = "'" + '1.000'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'" + '  -0.04%  '
$ws.Range('E6').Style = "Normal"

# Row 7
$ws.Range('D7').Value = "'" + '0.4461'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = "'" + '  +5.65%  '
$ws.Range('E7').Style = "Normal"

# Row 8
$ws.Range('E8').Value = "'" + '  +0.43%  '
$ws.Range('E8').Style = "Normal"

# Row 9
$ws.Range('D9').Value = "'" + '0.07345'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'" + '  +1.76%  '
$ws.Range('E9').Style = "Normal"

# Row 10
$ws.Range('D10').Value = "'" + '0.8593'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'" + '  +1.60%  '
$ws.Range('E10').Style = "Normal"

# Row 11
$ws.Range('D11').Value = "'" + '20.67'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'" + '  -1.07%  '
$ws.Range('E11').Style = "Normal"

# Row 12
$ws.Range('D12').Value = "'" + '1.798.12'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'" + '  -1.31%  '
$ws.Range('E12').Style = "Normal"

# Row 13
$ws.Range('D13').Value = "'" + '6.631'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'" + '  -0.13%  '
$ws.Range('E13').Style = "Normal"

# Row 14
$ws.Range('D14').Value = "'" + '92.35'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'" + '  +3.49%  '
$ws.Range('E14').Style = "Normal"

# Row 15
$ws.Range('D15').Value = "'" + '0.07064'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'" + '  -0.16%  '
$ws.Range('E15').Style = "Normal"

# Row 16
$ws.Range('D16').Value = "'" + '5.269'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'" + '  -0.29%  '
$ws.Range('E16').Style = "Normal"

# Row 17
$ws.Range('E17').Value = "'" + '  -0.01%  '
$ws.Range('E17').Style = "Normal"

# Row 18
$ws.Range('D18').Value = "'" + '0.000008688'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'" + '  -1.61%  '
$ws.Range('E18').Style = "Normal"

# Row 19
$ws.Range('E19').Value = "'" + '  -0.18%  '
$ws.Range('E19').Style = "Normal"

# Row 20
$ws.Range('E20').Value = "'" + '  -1.18%  '
$ws.Range('E20').Style = "Normal"

# Row 21
$ws.Range('D21').Value = "'" + '26.761.57'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'" + '  -1.40%  '
$ws.Range('E21').Style = "Normal"

# Row 22
$ws.Range('D22').Value = "'" + '5.162'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "'" + '  +0.92%  '
$ws.Range('E22').Style = "Normal"

# Row 23
$ws.Range('D23').Value = "'" + '10.79'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'" + '  -0.39%  '
$ws.Range('E23').Style = "Normal"

# Row 24
$ws.Range('D24').Value = "'" + '1.986'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'" + '  +0.67%  '
$ws.Range('E24').Style = "Normal"

# Row 25
$ws.Range('D25').Value = "'" + '151.82'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'" + '  -0.03%  '
$ws.Range('E25').Style = "Normal"

# Row 26
$ws.Range('D26').Value = "'" + '18.40'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'" + '  +0.38%  '
$ws.Range('E26').Style = "Normal"

# Row 27
$ws.Range('D27').Value = "'" + '2.171'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'" + '  -3.55%  '
$ws.Range('E27').Style = "Normal"

# Row 28
$ws.Range('D28').Value = "'" + '5.196'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'" + '  -0.26%  '
$ws.Range('E28').Style = "Normal"

# Row 29
$ws.Range('D29').Value = "'" + '117.34'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'" + '  +1.05%  '
$ws.Range('E29').Style = "Normal"

# Row 30
$ws.Range('E30').Value = "'" + '  -0.03%  '
$ws.Range('E30').Style = "Normal"

# Row 31
$ws.Range('D31').Value = "'" + '0.7409'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'" + '  -0.03%  '
$ws.Range('E31').Style = "Normal"

# Row 32
$ws.Range('D32').Value = "'" + '1.157'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'" + '  -1.96%  '
$ws.Range('E32').Style = "Normal"

# Row 33
$ws.Range('B33').Value = "'" + 'HuobiToken'
$ws.Range('B33').Style = "Normal"
$ws.Range('C33').Value = "'" + 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C33').Style = "Normal"
$ws.Range('D33').Value = "'" + '2.917'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'" + '  -1.42%  '
$ws.Range('E33').Style = "Normal"

# Row 34
$ws.Range('B34').Value = "'" + 'Filecoin'
$ws.Range('B34').Style = "Normal"
$ws.Range('C34').Value = "'" + 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C34').Style = "Normal"
$ws.Range('D34').Value = "'" + '4.453'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'" + '  +0.71%  '
$ws.Range('E34').Style = "Normal"

# Row 35
$ws.Range('D35').Value = "'" + '1.000'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'" + '  -0.03%  '
$ws.Range('E35').Style = "Normal"

# Row 36
$ws.Range('D36').Value = "'" + '1.083'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "'" + '  -1.39%  '
$ws.Range('E36').Style = "Normal"

# Row 37
$ws.Range('D37').Value = "'" + '0.01958'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'" + '  -0.52%  '
$ws.Range('E37').Style = "Normal"

# Row 38
$ws.Range('D38').Value = "'" + '0.05192'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'" + '  -0.96%  '
$ws.Range('E38').Style = "Normal"

# Row 39
$ws.Range('D39').Value = "'" + '0.5287'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'" + '  +5.05%  '
$ws.Range('E39').Style = "Normal"

# Row 40
$ws.Range('E40').Value = "'" + '  -1.36%  '
$ws.Range('E40').Style = "Normal"

# Row 41
$ws.Range('D41').Value = "'" + '6.966'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "'" + '  -4.57%  '
$ws.Range('E41').Style = "Normal"

# Row 42
$ws.Range('D42').Value = "'" + '0.1682'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'" + '  -0.69%  '
$ws.Range('E42').Style = "Normal"

# Row 43
$ws.Range('E43').Value = "'" + '  +7.08%  '
$ws.Range('E43').Style = "Normal"

# Row 44
$ws.Range('D44').Value = "'" + '8.433'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'" + '  -1.80%  '
$ws.Range('E44').Style = "Normal"

# Row 45
$ws.Range('B45').Value = "'" + 'RenderToken'
$ws.Range('B45').Style = "Normal"
$ws.Range('C45').Value = "'" + 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C45').Style = "Normal"
$ws.Range('D45').Value = "'" + '1.978'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'" + '  +4.62%  '
$ws.Range('E45').Style = "Normal"

# Row 46
$ws.Range('B46').Value = "'" + 'EnergySwap'
$ws.Range('B46').Style = "Normal"
$ws.Range('C46').Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C46').Style = "Normal"
$ws.Range('D46').Value = "'" + '10.43'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'" + '  -1.51%  '
$ws.Range('E46').Style = "Normal"

# Row 47
$ws.Range('D47').Value = "'" + '105.03'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'" + '  -0.98%  '
$ws.Range('E47').Style = "Normal"

# Row 48
$ws.Range('D48').Value = "'" + '1.672'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'" + '  +1.15%  '
$ws.Range('E48').Style = "Normal"

# Row 49
$ws.Range('D49').Value = "'" + '0.9996'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'" + '  -0.12%  '
$ws.Range('E49').Style = "Normal"

# Row 50
$ws.Range('D50').Value = "'" + '0.06293'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'" + '  -1.22%  '
$ws.Range('E50').Style = "Normal"

# Row 51
$ws.Range('D51').Value = "'" + '0.9164'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'" + '  +1.24%  '
$ws.Range('E51').Style = "Normal"
